$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 151 (shifts old rows 151-155 down to 152-156) ---
$ws.Rows.Item(151).Insert()

# Fix up the style of the new A151 cell (bordered/centered header-row style) by
# copying the format from the cell below (which carries the same style).
$ws.Cells.Item(152, 1).Copy()
$ws.Cells.Item(151, 1).PasteSpecial(-4122)  # xlPasteFormats

# --- Populate the brand-new row 151 with the new match data ---
$ws.Cells.Item(151, 1).Value2  = 149
$ws.Cells.Item(151, 2).Value2  = 7993785
$ws.Cells.Item(151, 3).Value2  = "Croatia HNL"
$ws.Cells.Item(151, 4).Value2  = "Croatia HNL"
$ws.Cells.Item(151, 5).Value2  = 45399.5
$ws.Cells.Item(151, 6).Value2  = "Dinamo Zagreb"
$ws.Cells.Item(151, 7).Value2  = "NK Varazdin"
$ws.Cells.Item(151, 8).Value2  = 2
$ws.Cells.Item(151, 9).Value2  = 1
$ws.Cells.Item(151, 10).Value2 = "H"
$ws.Cells.Item(151, 11).Value2 = 1.2
$ws.Cells.Item(151, 12).Value2 = 6.5
$ws.Cells.Item(151, 13).Value2 = 13
$ws.Cells.Item(151, 14).Value2 = 1.181
$ws.Cells.Item(151, 15).Value2 = 6.5
$ws.Cells.Item(151, 16).Value2 = 15
$ws.Cells.Item(151, 17).Value2 = -2
$ws.Cells.Item(151, 18).Value2 = 1.925
$ws.Cells.Item(151, 19).Value2 = 1.925
$ws.Cells.Item(151, 20).Value2 = 3.25
$ws.Cells.Item(151, 21).Value2 = 2.025
$ws.Cells.Item(151, 22).Value2 = 1.825
$ws.Cells.Item(151, 23).Value2 = 0.181
$ws.Cells.Item(151, 24).Value2 = -1
$ws.Cells.Item(151, 25).Value2 = -1
$ws.Cells.Item(151, 26).Value2 = -1
$ws.Cells.Item(151, 27).Value2 = 0.925
$ws.Cells.Item(151, 28).Value2 = -0.5
$ws.Cells.Item(151, 29).Value2 = 0.4125

# --- Renumber the sequential "id" column (A) for the rows that got shifted down ---
$ws.Cells.Item(152, 1).Value2 = 150
$ws.Cells.Item(153, 1).Value2 = 151
$ws.Cells.Item(154, 1).Value2 = 152
$ws.Cells.Item(155, 1).Value2 = 153
$ws.Cells.Item(156, 1).Value2 = 154

# --- Update a handful of closing-odds figures that were revised for the two
#     matches now sitting in rows 152 and 153 (previously rows 151 and 152) ---
$ws.Cells.Item(152, 15).Value2 = 3.6    # O152
$ws.Cells.Item(152, 16).Value2 = 6.5    # P152

$ws.Cells.Item(153, 14).Value2 = 5.25   # N153
$ws.Cells.Item(153, 15).Value2 = 3.8    # O153
$ws.Cells.Item(153, 16).Value2 = 1.615  # P153
$ws.Cells.Item(153, 18).Value2 = 2.025  # R153
$ws.Cells.Item(153, 19).Value2 = 1.825  # S153
$ws.Cells.Item(153, 21).Value2 = 1.95   # U153
$ws.Cells.Item(153, 22).Value2 = 1.9    # V153
